$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date as an Excel serial
# number. Rows 2-39 all moved from serial 45188 (2023-09-19) to serial
# 45189 (2023-09-20).
for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
